$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-separated names that should use periods (typo/formatting fix) ---
$nameFixes = @(
    @{Addr='E45'; Val='MONROY. AGUSTIN ALEJANDRO'}
    @{Addr='E57'; Val='MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'}
    @{Addr='E77'; Val='BOFFELLI. MARIA INES'}
    @{Addr='E130'; Val='ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'}
    @{Addr='E141'; Val='SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'}
)
foreach ($item in $nameFixes) {
    $ws.Range($item.Addr).Value = $item.Val
}

# --- Fix Importe (amount) column: remove thousand separators, use "." as decimal separator ---
$amountFixes = @(
    @{Addr='H2'; Val='1800.00'}
    @{Addr='H3'; Val='13480.00'}
    @{Addr='H4'; Val='5300.00'}
    @{Addr='H5'; Val='272000.00'}
    @{Addr='H6'; Val='300.00'}
    @{Addr='H7'; Val='1530.40'}
    @{Addr='H8'; Val='154000.00'}
    @{Addr='H9'; Val='6874.00'}
    @{Addr='H10'; Val='8608.88'}
    @{Addr='H11'; Val='313685.39'}
    @{Addr='H12'; Val='249386.73'}
    @{Addr='H13'; Val='30066.56'}
    @{Addr='H14'; Val='19.00'}
    @{Addr='H15'; Val='154729.23'}
    @{Addr='H16'; Val='10162.50'}
    @{Addr='H17'; Val='18791.70'}
    @{Addr='H18'; Val='9539.00'}
    @{Addr='H19'; Val='13834.23'}
    @{Addr='H20'; Val='4822.00'}
    @{Addr='H21'; Val='21773.01'}
    @{Addr='H22'; Val='1000.00'}
    @{Addr='H23'; Val='300.00'}
    @{Addr='H24'; Val='553.00'}
    @{Addr='H25'; Val='11273.33'}
    @{Addr='H26'; Val='1070.00'}
    @{Addr='H27'; Val='3931.37'}
    @{Addr='H28'; Val='1781.29'}
    @{Addr='H29'; Val='1589.38'}
    @{Addr='H30'; Val='696.00'}
    @{Addr='H31'; Val='11698.59'}
    @{Addr='H32'; Val='8620.20'}
    @{Addr='H33'; Val='1760.00'}
    @{Addr='H34'; Val='4074.20'}
    @{Addr='H35'; Val='27.90'}
    @{Addr='H36'; Val='12660.00'}
    @{Addr='H37'; Val='7902.02'}
    @{Addr='H38'; Val='1258.30'}
    @{Addr='H39'; Val='4342.00'}
    @{Addr='H40'; Val='2695.00'}
    @{Addr='H41'; Val='654.70'}
    @{Addr='H42'; Val='1133.50'}
    @{Addr='H43'; Val='1940.00'}
    @{Addr='H44'; Val='15615.00'}
    @{Addr='H45'; Val='3900.00'}
    @{Addr='H46'; Val='6849.00'}
    @{Addr='H47'; Val='810.00'}
    @{Addr='H48'; Val='219948.00'}
    @{Addr='H49'; Val='3488.95'}
    @{Addr='H50'; Val='1339.50'}
    @{Addr='H51'; Val='280.00'}
    @{Addr='H52'; Val='4982.00'}
    @{Addr='H53'; Val='490.00'}
    @{Addr='H54'; Val='14700.00'}
    @{Addr='H55'; Val='5832.00'}
    @{Addr='H56'; Val='17990.00'}
    @{Addr='H57'; Val='1100.00'}
    @{Addr='H58'; Val='2130.00'}
    @{Addr='H59'; Val='11400.00'}
    @{Addr='H60'; Val='1250.00'}
    @{Addr='H61'; Val='5227.00'}
    @{Addr='H62'; Val='6000.00'}
    @{Addr='H63'; Val='9483.54'}
    @{Addr='H64'; Val='14.10'}
    @{Addr='H65'; Val='415230.93'}
    @{Addr='H66'; Val='53261.49'}
    @{Addr='H67'; Val='33150.00'}
    @{Addr='H68'; Val='15.95'}
    @{Addr='H69'; Val='59.99'}
    @{Addr='H70'; Val='2880.00'}
    @{Addr='H71'; Val='1200.00'}
    @{Addr='H72'; Val='53.00'}
    @{Addr='H73'; Val='14836.00'}
    @{Addr='H74'; Val='29985.40'}
    @{Addr='H75'; Val='29696.00'}
    @{Addr='H76'; Val='26894.00'}
    @{Addr='H77'; Val='18012.00'}
    @{Addr='H78'; Val='110.00'}
    @{Addr='H79'; Val='11985.40'}
    @{Addr='H80'; Val='2427.00'}
    @{Addr='H81'; Val='9184.00'}
    @{Addr='H82'; Val='30898.22'}
    @{Addr='H83'; Val='7444.90'}
    @{Addr='H84'; Val='1296.00'}
    @{Addr='H85'; Val='611.76'}
    @{Addr='H86'; Val='3410.50'}
    @{Addr='H87'; Val='2226.60'}
    @{Addr='H88'; Val='350.00'}
    @{Addr='H89'; Val='229.24'}
    @{Addr='H90'; Val='160.00'}
    @{Addr='H91'; Val='1380.00'}
    @{Addr='H92'; Val='810.00'}
    @{Addr='H93'; Val='1600.00'}
    @{Addr='H94'; Val='3124.00'}
    @{Addr='H95'; Val='22650.00'}
    @{Addr='H96'; Val='50010.00'}
    @{Addr='H97'; Val='6840.00'}
    @{Addr='H98'; Val='700.00'}
    @{Addr='H99'; Val='3500.00'}
    @{Addr='H100'; Val='60748.24'}
    @{Addr='H101'; Val='973.72'}
    @{Addr='H102'; Val='73.25'}
    @{Addr='H103'; Val='7340.90'}
    @{Addr='H104'; Val='1260000.00'}
    @{Addr='H105'; Val='18447.60'}
    @{Addr='H106'; Val='8000.00'}
    @{Addr='H107'; Val='10000.00'}
    @{Addr='H108'; Val='4000.00'}
    @{Addr='H109'; Val='22264.00'}
    @{Addr='H110'; Val='2500.00'}
    @{Addr='H111'; Val='2762.50'}
    @{Addr='H112'; Val='2556.00'}
    @{Addr='H113'; Val='3000.00'}
    @{Addr='H114'; Val='2000.00'}
    @{Addr='H115'; Val='1500.00'}
    @{Addr='H116'; Val='4000.00'}
    @{Addr='H117'; Val='3800.00'}
    @{Addr='H118'; Val='19995.00'}
    @{Addr='H119'; Val='9820.00'}
    @{Addr='H120'; Val='6000.00'}
    @{Addr='H121'; Val='2500.00'}
    @{Addr='H122'; Val='1500.00'}
    @{Addr='H123'; Val='17880.00'}
    @{Addr='H124'; Val='6000.00'}
    @{Addr='H125'; Val='780.00'}
    @{Addr='H126'; Val='605.00'}
    @{Addr='H127'; Val='915.52'}
    @{Addr='H128'; Val='2590.00'}
    @{Addr='H129'; Val='1780.00'}
    @{Addr='H130'; Val='7240.00'}
    @{Addr='H131'; Val='12295.00'}
    @{Addr='H132'; Val='2369.65'}
    @{Addr='H133'; Val='12317.00'}
    @{Addr='H134'; Val='3720.00'}
    @{Addr='H135'; Val='18000.00'}
    @{Addr='H136'; Val='2994.16'}
    @{Addr='H137'; Val='690.00'}
    @{Addr='H138'; Val='2100.00'}
    @{Addr='H139'; Val='27686.26'}
    @{Addr='H140'; Val='810.00'}
    @{Addr='H141'; Val='2300.00'}
    @{Addr='H142'; Val='3500.00'}
    @{Addr='H143'; Val='5500.00'}
    @{Addr='H144'; Val='2772.35'}
    @{Addr='H145'; Val='4323.19'}
    @{Addr='H146'; Val='1304.40'}
    @{Addr='H147'; Val='3920.00'}
    @{Addr='H148'; Val='1150.00'}
    @{Addr='H149'; Val='43200.00'}
    @{Addr='H150'; Val='900.00'}
    @{Addr='H151'; Val='1852045.68'}
    @{Addr='H152'; Val='1670.00'}
    @{Addr='H153'; Val='260000.00'}
    @{Addr='H154'; Val='4200.00'}
    @{Addr='H155'; Val='105000.00'}
    @{Addr='H156'; Val='105000.00'}
    @{Addr='H157'; Val='105000.00'}
    @{Addr='H158'; Val='175000.00'}
    @{Addr='H159'; Val='175000.00'}
    @{Addr='H160'; Val='245000.00'}
    @{Addr='H161'; Val='105000.00'}
    @{Addr='H162'; Val='105000.00'}
    @{Addr='H163'; Val='105000.00'}
    @{Addr='H164'; Val='105000.00'}
    @{Addr='H165'; Val='175000.00'}
    @{Addr='H166'; Val='315000.00'}
    @{Addr='H167'; Val='175000.00'}
    @{Addr='H168'; Val='105000.00'}
    @{Addr='H169'; Val='130000.00'}
    @{Addr='H170'; Val='105000.00'}
    @{Addr='H171'; Val='105000.00'}
    @{Addr='H172'; Val='105000.00'}
    @{Addr='H173'; Val='420000.00'}
    @{Addr='H174'; Val='5100.00'}
    @{Addr='H175'; Val='7500.00'}
    @{Addr='H176'; Val='9322.89'}
    @{Addr='H177'; Val='13160.00'}
    @{Addr='H178'; Val='18350.00'}
    @{Addr='H179'; Val='1800.00'}
    @{Addr='H180'; Val='47150.00'}
    @{Addr='H181'; Val='1152.00'}
)
foreach ($item in $amountFixes) {
    $r = $ws.Range($item.Addr)
    $r.NumberFormat = "@"
    $r.Value = $item.Val
    $r.Style = "Normal"
}
